$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column F (ECE2213, CMPE2213, EE2213) content, without shifting other cells
$ws.Range("F1:F3").ClearContents()

# Add new column H with Statistics / Transcript Reader related values
$ws.Range("H1").Value = "a"
$ws.Range("H2").Value = "b"
$ws.Range("H3").Value = "c"
